# Update "Training Dashboard" sheet with new progress as of 04-Nov-2025.
# Column H = PERIOD TO EXPIRE (decrements by 1 day), Column I = LAST UPDATE date.
# NOTE: this runtime's Range/Cells ".Value" getter is unreliable for reading
# back a value (it returns a member-signature placeholder rather than the
# actual contents), so ".Value2" is used for reads. Writing plain date-like
# strings through ".Value" causes automatic conversion to a date serial
# number, so a leading apostrophe is used to force the LAST UPDATE column to
# stay a text value (matching the original inline-string cells), which Excel
# strips from the stored text while only flipping the cell's quotePrefix
# flag (style/format stay the same).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

$newLastUpdate = "04-Nov-2025"

foreach ($r in 3..9) {
    $hCell = $ws.Cells.Item($r, 8)   # column H: PERIOD TO EXPIRE
    $hCell.Value = $hCell.Value2 - 1

    $iCell = $ws.Cells.Item($r, 9)   # column I: LAST UPDATE
    $iCell.Value = "'" + $newLastUpdate
}
